# Automatic update of files.
# Applies the record-level corrections to the "Artfynd" sheet:
#  - Rows 2/3 swap their Id/Ost/Nord/Starttid/Sluttid (two duplicate
#    "Garnlav" observations had their coordinates & times mixed up).
#  - Rows 6/7 swap their whole content (a "Garnlav" row and a
#    "Tretåig hackspett" row had been entered on the wrong lines).
#  - Rows 25/26 swap their whole content in the same way.
#  - Several rows get their Taxonsorteringsordning (column B) bumped
#    by one (a taxon-sort-order renumbering upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell {
    param($addr, $value)
    $ws.Range($addr).Value = $value
}

function Clear-Cell {
    param($addr)
    $ws.Range($addr).ClearContents()
}

# ---- Rows 2 & 3: swap Id / Ost / Nord / Starttid / Sluttid ----
Set-Cell "A2" 131046824
Set-Cell "Q2" 401653
Set-Cell "R2" 6818054
Set-Cell "Z2" "14:50"
Set-Cell "AB2" "14:50"

Set-Cell "A3" 131046825
Set-Cell "Q3" 401650
Set-Cell "R3" 6818017
Set-Cell "Z3" "14:52"
Set-Cell "AB3" "14:52"

# ---- Rows 6 & 7: full content swap ----
Set-Cell "A6" 131046830
Set-Cell "E6" 6425
Set-Cell "F6" "Garnlav"
Set-Cell "G6" "Alectoria sarmentosa"
Set-Cell "H6" "(Ach.) Ach."
Clear-Cell "M6"
Set-Cell "Q6" 401538
Set-Cell "R6" 6818009
Set-Cell "Z6" "15:15"
Set-Cell "AB6" "15:15"
Clear-Cell "AC6"

Set-Cell "A7" 131046772
Set-Cell "B7" 57884
Set-Cell "E7" 100109
Set-Cell "F7" "Tretåig hackspett"
Set-Cell "G7" "Picoides tridactylus"
Set-Cell "H7" "(Linnaeus, 1758)"
Set-Cell "M7" "färska spår"
Set-Cell "Q7" 401507
Set-Cell "R7" 6818011
Set-Cell "Z7" "15:17"
Set-Cell "AB7" "15:17"
Set-Cell "AC7" "Färska ringhack (tall)"

# ---- Rows 25 & 26: full content swap ----
Set-Cell "A25" 131047014
Set-Cell "B25" 57884
Set-Cell "E25" 100109
Set-Cell "F25" "Tretåig hackspett"
Set-Cell "G25" "Picoides tridactylus"
Set-Cell "H25" "(Linnaeus, 1758)"
Set-Cell "M25" "färska spår"
Set-Cell "Q25" 401378
Set-Cell "R25" 6818082
Set-Cell "Z25" "15:21"
Set-Cell "AB25" "15:21"
Set-Cell "AC25" "Troliga spår efter tretåig hackspett (barkfälkning)"
Set-Cell "AE25" $true

Set-Cell "A26" 131046832
Set-Cell "E26" 6425
Set-Cell "F26" "Garnlav"
Set-Cell "G26" "Alectoria sarmentosa"
Set-Cell "H26" "(Ach.) Ach."
Clear-Cell "M26"
Set-Cell "Q26" 401350
Set-Cell "R26" 6818162
Set-Cell "Z26" "15:24"
Set-Cell "AB26" "15:24"
Clear-Cell "AC26"
Set-Cell "AE26" $false

# ---- Column B (Taxonsorteringsordning) bumped by one ----
# Garnlav rows: 79243 -> 79244
foreach ($r in 2,3,4,8,10,12,18,23,24,27) {
    Set-Cell ("B" + $r) 79244
}
# The Garnlav row that landed on 6 and 26 after the swap above
Set-Cell "B6" 79244
Set-Cell "B26" 79244

# Vitgrynig nållav rows: 83223 -> 83224
foreach ($r in 14,17) {
    Set-Cell ("B" + $r) 83224
}

# Violettgrå tagellav row: 79275 -> 79276
Set-Cell "B16" 79276

# Blanksvart spiklav row: 78646 -> 78647
Set-Cell "B20" 78647

# Liten svartspik row: 78255 -> 78256
Set-Cell "B22" 78256
